$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.936.72'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '2.618.72'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '533.06'
$ws.Range('E5').Value = '  -1.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.46'
$ws.Range('E6').Value = '  +1.42%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.568'
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.93'
$ws.Range('E9').Value = '  +7.66%  '
$ws.Range('E10').Value = '  -1.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.335'
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').Value = '3.087.02'
$ws.Range('E13').Value = '  +0.87%  '
$ws.Range('D14').Value = '58.903.47'
$ws.Range('E14').Value = '  -0.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.93'
$ws.Range('E15').Value = '  +1.55%  '
$ws.Range('D16').Value = '2.626.27'
$ws.Range('E16').Value = '  +0.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000133'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.41'
$ws.Range('E18').Value = '  +1.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '336.16'
$ws.Range('E19').Value = '  -1.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.17'
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.20'
$ws.Range('E21').Value = '  -3.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.24'
$ws.Range('E23').Value = '  -2.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.414'
$ws.Range('E24').Value = '  +1.27%  '
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.20'
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('D28').Value = '0.0₃0737'
$ws.Range('E28').Value = '  -0.96%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  -2.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.82'
$ws.Range('E31').Value = '  +0.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '150.69'
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.72'
$ws.Range('E33').Value = '  -0.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.94'
$ws.Range('E34').Value = '  -0.88%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.12'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.832'
$ws.Range('E36').Value = '  +0.79%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.824'
$ws.Range('E37').Value = '  -1.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.43'
$ws.Range('E38').Value = '  -2.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.57'
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '282.73'
$ws.Range('E40').Value = '  +2.80%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.594'
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.73'
$ws.Range('E43').Value = '  -0.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0534'
$ws.Range('E44').Value = '  +2.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.95'
$ws.Range('E45').Value = '  +1.81%  '
$ws.Range('E46').Value = '  -1.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0224'
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('D48').Value = '1.944.76'
$ws.Range('E48').Value = '  -0.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.51'
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.06'
$ws.Range('E50').Value = '  -2.47%  '
$ws.Range('E51').Value = '  +0.39%  '
